# Apply the edits described by the diff:
# 1. Update the "Cloud - C (KRTA9AA3)" course string to include the second group code.
# 2. Update the bare "KRTA9AA3" code string to include the second group code.
# 3. Fill in room values in column F for rows 3,4,7,8,11,12,15,16:
#      rows 3,4,11,12 -> "U3-Amphi"
#      rows 7,8,15,16 -> "U3-4"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1 & 2: update shared text that appears in columns A and B ---
$rows = @(3, 4, 7, 8, 11, 12, 15, 16)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 1).Value = "Cloud - C (KRTA9AA3/KUPT9BB1)"
    $ws.Cells.Item($r, 2).Value = "KRTA9AA3/KUPT9BB1"
}

# --- 3: fill column F (room) ---
$u3Amphi = @(3, 4, 11, 12)
foreach ($r in $u3Amphi) {
    $ws.Cells.Item($r, 6).Value = "U3-Amphi"
}

$u34 = @(7, 8, 15, 16)
foreach ($r in $u34) {
    $ws.Cells.Item($r, 6).Value = "U3-4"
}
